# Doing Updates for Financials
# Update the historical (2017 / column E) figures on the ETCC sheet, plus the
# 5 figures that changed across multiple years for "Capital Expenditures" (row 91).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ETCC")

# Balance Sheet values - column E (period ending 2017) updates
$ws.Range("E41").Value = 300      # Cash And Cash Equivalents
$ws.Range("E43").Value = 5900     # Net Receivables
$ws.Range("E44").Value = 2400     # Inventory
$ws.Range("E45").Value = 24900    # Other Current Assets
$ws.Range("E46").Value = 33500    # Total Current Assets
$ws.Range("E48").Value = 12400    # Property Plant and Equipment
$ws.Range("E52").Value = 500      # Other Assets
$ws.Range("E57").Value = 3900     # Accounts Payable
$ws.Range("E58").Value = 0        # Short/Current Long Term Debt
$ws.Range("E59").Value = 22100    # Other Current Liabilities
$ws.Range("E60").Value = 15200    # Total Current Liabilities
$ws.Range("E61").Value = 20900    # Long Term Debt
$ws.Range("E66").Value = 36900    # Total Liabilities
$ws.Range("E72").Value = -10700   # Retained Earnings
$ws.Range("E76").Value = 9500     # Total Stockholder Equity

# Cash Flow Statement - Capital Expenditures (row 91) updates across several years
$ws.Range("D91").Value = -400
$ws.Range("E91").Value = -600
$ws.Range("F91").Value = -1100
$ws.Range("H91").Value = -1300
$ws.Range("J91").Value = -1900
